# Update the fiscal-year 2081/82 row (row 12) on the sheet with the
# "complete fiscal year" figures per the commit message.
# A leading apostrophe forces Excel to store each value as literal text
# (matching the existing inlineStr/text cells in this row) instead of
# auto-converting the numeric-looking strings to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "' 1804122731 "
$ws.Range("C12").Value = "' 277030202 "
$ws.Range("D12").Value = "' 1527092530 "
$ws.Range("E12").Value = "' 2081152933 "
$ws.Range("F12").Value = "'6.51"
$ws.Range("G12").Value = "'13.31"
$ws.Range("H12").Value = "'86.69"
